# add area to Q files stn4
#
# Adds an "Area" / "Atotal" computation (columns G/H) alongside the
# existing "segment" / "Q" / "Qtotal" computation (columns D/E/F), plus a
# small two-column summary (J/K) that mirrors the grand totals. Also fills
# in the previously blank/placeholder depth (B) and velocity (C) values for
# the first and last verticals (0 instead of a "-" placeholder), matching
# the fully-computed discharge sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers (row 1) -----------------------------------------------
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Fill in the edge verticals that used to hold a "-" placeholder -----
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0

# --- Area-per-segment column (G), mirrors the Q-per-segment column (E) -
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4").Formula = "=(D4-D3)*B4/100"
$ws.Range("G5").Formula = "=(D5-D4)*B5/100"
$ws.Range("G6").Formula = "=(D6-D5)*B6/100"
$ws.Range("G7").Formula = "=(D7-D6)*B7/100"
$ws.Range("G8").Formula = "=(D8-D7)*B8/100"
$ws.Range("G9").Formula = "=(D9-D8)*B9/100"
$ws.Range("G10").Formula = "=(D10-D9)*B10/100"
$ws.Range("G11").Formula = "=(D11-D10)*B11/100"
$ws.Range("G12").Formula = "=(D12-D11)*B12/100"
$ws.Range("G13").Formula = "=(D13-D12)*B13/100"
$ws.Range("G14").Formula = "=(D14-D13)*B14/100"
$ws.Range("G15").Formula = "=(D15-D14)*B15/100"

# --- Grand total of the area column (H), mirrors Qtotal (F) ------------
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Small summary block (J2:K2) mirroring the two grand totals --------
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Selection / view, matching the authored state ----------------------
[void]$ws.Range("J2:K2").Select()
